# Weekly CompStat data refresh: new report week (Volume/Number + date range)
# and updated crime-complaint figures for rows 16-31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump the report Volume/Number and the covered week dates ---
$titleCell = $ws.Cells.Item(8, 1)   # "Volume 31   Number  39" -> "...  40"
$titleCell.Characters(21, 2).Text = "40"

$weekCell = $ws.Cells.Item(9, 3)    # "Report Covering the Week  9/23/2024  Through  9/29/2024"
$weekCell.Characters(27, 9).Text = "9/30/2024"
$weekCell.Characters(47, 9).Text = "10/6/2024"

# --- Crime-complaint table: cells whose value flips between a number and the
#     "no activity" text markers ("0" / "***.*") need both their content and
#     their style swapped to match; use a donor cell with the right style already
#     applied and PasteSpecial so no brand-new style entries get created. ---
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C24").PasteSpecial(-4122)

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 1

$ws.Range("C16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D24").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 2

$ws.Range("C16").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 2

$ws.Range("L14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = -100

$ws.Range("L14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -100

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E24").PasteSpecial(-4122)

$ws.Range("L14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100

$ws.Range("L14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G16").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G17").Value = 1

$ws.Range("C16").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = 2

$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H16").PasteSpecial(-4122)

$ws.Range("L14").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("H17").Value = -100

$ws.Range("L14").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = -100

$excel.CutCopyMode = $false

# --- Remaining cells: value-only updates (style/type unchanged) ---
$ws.Range("D19").Value = 4
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = -83.333333333333
$ws.Range("F16").Value = 3
$ws.Range("F21").Value = 5
$ws.Range("F24").Value = 3
$ws.Range("G18").Value = 2
$ws.Range("G19").Value = 7
$ws.Range("G21").Value = 10
$ws.Range("G24").Value = 3
$ws.Range("G26").Value = 7
$ws.Range("G28").Value = 1
$ws.Range("H19").Value = -71.428571428571
$ws.Range("H21").Value = -50
$ws.Range("H24").Value = 0
$ws.Range("H26").Value = -42.857142857142
$ws.Range("H28").Value = 200
$ws.Range("I16").Value = 37
$ws.Range("I21").Value = 87
$ws.Range("J17").Value = 8
$ws.Range("J18").Value = 5
$ws.Range("J19").Value = 44
$ws.Range("J21").Value = 75
$ws.Range("J26").Value = 43
$ws.Range("J31").Value = 4
$ws.Range("K16").Value = 117.647058823529
$ws.Range("K17").Value = 37.5
$ws.Range("K18").Value = -80
$ws.Range("K19").Value = -18.181818181818
$ws.Range("K21").Value = 16
$ws.Range("K26").Value = -27.906976744186
$ws.Range("K31").Value = -50
$ws.Range("L16").Value = 76.190476190476
$ws.Range("L17").Value = -15.384615384615
$ws.Range("L21").Value = 50
$ws.Range("L26").Value = 19.230769230769
$ws.Range("M16").Value = 68.181818181818
$ws.Range("M19").Value = -38.983050847457
$ws.Range("M21").Value = -6.451612903225
$ws.Range("M24").Value = -52.941176470588
$ws.Range("N16").Value = -77.976190476190
$ws.Range("N17").Value = -67.647058823529
$ws.Range("N19").Value = -75.510204081632
$ws.Range("N21").Value = -77.862595419847
